# Updated cryptos list on Fri Nov  8 17:46:10 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) columns of
# the crypto-price table on rows 2-51 of the active sheet with freshly
# scraped values. Both columns are plain text cells in the workbook (prices
# use '.' as a thousands separator in places, and the volume column carries
# literal leading/trailing padding spaces around the percentage), so every
# write goes through .Value as a string.
#
# A few of the new price strings (e.g. "1.00", "198.80") are, character for
# character, things Excel's COM layer would normally auto-coerce into a
# Number when assigned via Range.Value. To keep those cells as text (matching
# the rest of column D and the original file's inlineStr cells), the format
# is forced to "@" (Text) immediately before the write and then reset back to
# the default "Normal" style afterwards so no stray per-cell formatting is
# left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.236.75"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "2.912.01"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "591.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.546"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("D10").Value = "2.910.41"
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("E11").Value = "  +14.17%  "
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").Value = "3.443.86"
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("D15").Value = "76.069.19"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000186"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "2.905.85"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.00%  "
$ws.Range("E20").Value = "  -5.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "368.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.75%  "
$ws.Range("E22").Value = "  +3.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").Value = "3.045.75"
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.14%  "
$ws.Range("E32").Value = "  -4.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "493.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.31%  "
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.386"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.108"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +19.19%  "
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.109"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "177.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.57%  "
$ws.Range("E45").Value = "  -3.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("E47").Value = "  -5.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.579"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.37%  "
